# Update the "Last Updated" timestamp on the Metadata sheet.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value2 = "05 Nov 2025, 11:00 AM"

# Stock List sheet: a new row (CAPTRU-RE1) was inserted at row 2, pushing
# every existing data row down by one and dropping the former last row
# (row 76 / TRAVELFOOD) off the bottom of the table.
$ws = $wb.Worksheets.Item("Stock List")

for ($i = 76; $i -ge 3; $i--) {
    $src = $i - 1
    $ws.Cells.Item($i, 2).Value2 = $ws.Cells.Item($src, 2).Value2
    $ws.Cells.Item($i, 3).Value2 = $ws.Cells.Item($src, 3).Value2
    $ws.Cells.Item($i, 4).Value2 = $ws.Cells.Item($src, 4).Value2
    $ws.Cells.Item($i, 5).Value2 = $ws.Cells.Item($src, 5).Value2
    $ws.Cells.Item($i, 8).Value2 = $ws.Cells.Item($src, 8).Value2
}

$ws.Cells.Item(2, 2).Value2 = "CAPTRU-RE1"
$ws.Cells.Item(2, 3).Value2 = "CAPTRU-RE1"
$ws.Cells.Item(2, 4).Value2 = 5.67
$ws.Cells.Item(2, 5).Value2 = -11.9565
